# Update column F (rows 2-25) on Sheet1 with the new computed values
# for the "case with 380 kV" results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0.4443680307746121
    3  = 0.3878228170618172
    4  = 0.3531389305169483
    5  = 0.3390132514313251
    6  = 0.336668177824194
    7  = 0.3529483938345521
    8  = 0.4248636149814047
    9  = 0.5661985755041457
    10 = 0.6702781546542269
    11 = 0.7176906081379002
    12 = 0.7356546913071611
    13 = 0.7317853510981394
    14 = 0.7191683204515442
    15 = 0.7114413442032514
    16 = 0.6671810134426437
    17 = 0.6400460337125793
    18 = 0.6244449056556647
    19 = 0.619163680173358
    20 = 0.642933953830422
    21 = 0.7228739723492197
    22 = 0.7751780083420101
    23 = 0.7472568307830727
    24 = 0.6416283278901602
    25 = 0.5279251897347166
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
